# Applies the "Tạo xong bàn cờ" commit:
#  - Update row 5 (the "Vẽ bàn cờ" task row): the description text changes
#    from the old "16*16" wording to the new "19x23" wording, and the
#    actual start/end dates (08 tháng 10) are now filled in for H5/I5.
#  - Update the current selection / scroll position to reflect where the
#    author was working (C2 in view, I5 selected) when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# G5: updated description of the "create board" task
$ws.Range("G5").Value = "Tạo ra bàn cờ cơ bàn, xây dựng thuật toán tạo bàn cờ 19x23  "

# H5 / I5: actual start & end date now recorded as "08 tháng 10"
$ws.Range("H5").Value = "08 tháng 10"
$ws.Range("I5").Value = "08 tháng 10"

# Reflect the new view/selection state saved with the workbook
$ws.Range("C2").Select()
$ws.Range("I5").Select()
